$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmulatorData")

# Rename headers: "App ID" -> "AppID", "Tax Amount" -> "TaxAmount"
$ws.Range("C1").Value = "AppID"
$ws.Range("F1").Value = "TaxAmount"

# Row 2: shorten the note text
$ws.Range("A2").Value = "PayNow NoCF"

# Row 3: clear the Tax Amount value (F3)
$ws.Range("F3").ClearContents()

# Row 10: rename note text and clear the Tax Amount value (F10)
$ws.Range("A10").Value = "Populate only req fields NoCF"
$ws.Range("F10").ClearContents()

# Column A width widened (auto-fit after the longer row 9 text)
$ws.Columns.Item(1).ColumnWidth = 44.6

# Restore the last active selection to E13 like the saved file
$ws.Range("E13").Select()
